# Sprint 4 burndown sheet update:
#  - "Esperado" (expected) start value B5 goes from 26 -> 28 story points,
#    which ripples through the shared formula in B6:B11 (via E5 = B5/6).
#  - "Real" (actual) value for day 2 (C6) goes from 16 -> 10.
#  - C5 (actual) is kept in sync with the new starting point, 26 -> 28.
#  - Leave the cursor parked on M23, matching where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 4")

$ws.Range("B5").Value = 28
$ws.Range("C5").Value = 28
$ws.Range("C6").Value = 10

$ws.Range("M23").Select() | Out-Null
